$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Anzac Day (additional)" for 2026 right after the 2026 Anzac Day row ---
# Original row 9 (2026-06-08 / King's Birthday) gets pushed down; new row 9 is inserted.
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "'2026-04-27"
$ws.Range("A9").ClearFormats()
$ws.Range("B9").Value = "Anzac Day (additional)"

# --- Insert "Anzac Day (additional)" for 2027 right after the 2027 Anzac Day row ---
# After the first insertion, the 2027 Anzac Day row (originally row 21) is now row 22.
# Insert a fresh row at 23 for the new observance.
$ws.Rows("23:23").Insert()
$ws.Range("A23").Value = "'2027-04-26"
$ws.Range("A23").ClearFormats()
$ws.Range("B23").Value = "Anzac Day (additional)"
